# "Need to order.xlsx" - add encoder and phototransistor parts (with
# purchase links) that were shared in Discord, as per the commit message.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Data rows first (mirrors the order the info was actually typed in) ---

# Row 2: Encoder
$ws.Range("A2").Value = "Encoder"
$ws.Range("B2").Value = 2
$ws.Range("C2").Value = "https://www.mouser.com/ProductDetail/Bourns/PAC18R1-43D19F?qs=IS%252B4QmGtzzoiPnRwlkZJXA%3D%3D"
$ws.Hyperlinks.Add($ws.Range("C2"), "https://www.mouser.com/ProductDetail/Bourns/PAC18R1-43D19F?qs=IS%252B4QmGtzzoiPnRwlkZJXA%3D%3D", "", "https://www.mouser.com/ProductDetail/Bourns/PAC18R1-43D19F?qs=IS%252B4QmGtzzoiPnRwlkZJXA%3D%3D") | Out-Null

# Row 3: Phototransistors (name first, link added later below)
$ws.Range("A3").Value = "Phototransistors"
$ws.Range("B3").Value = 8

# --- Header row, added after the data rows ---
$ws.Range("A1").Value = "Part"
$ws.Range("B1").Value = "Quantity"
$ws.Range("C1").Value = "Link"
$ws.Rows(1).Font.Bold = $true

# --- Finish the Phototransistors row with its link ---
$ws.Range("C3").Value = "https://www.digikey.com/en/products/detail/vishay-semiconductor-opto-division/TEPT5700/1681193"
$ws.Hyperlinks.Add($ws.Range("C3"), "https://www.digikey.com/en/products/detail/vishay-semiconductor-opto-division/TEPT5700/1681193", "", "https://www.digikey.com/en/products/detail/vishay-semiconductor-opto-division/TEPT5700/1681193") | Out-Null

# --- Cosmetic touches to match the saved workbook ---
$ws.Columns("A").ColumnWidth = 13.2
$ws.PageSetup.Orientation = 1
$ws.Range("G8").Select() | Out-Null
